# Update column F ("dSF") values for rows 2-26 (repull data / push all data / mean calculation)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    3  = -7
    4  = -3
    5  = 1
    6  = 2
    7  = -2
    9  = -2
    10 = -2
    11 = -1
    12 = -4
    14 = -3
    15 = 3
    16 = 2
    17 = -4
    18 = 2
    19 = -1
    20 = -4
    21 = -1
    22 = 3
    23 = 9
    24 = 3
    26 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
